$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 205, pushing existing rows 205-276 down to 206-277.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new record's data.
$ws.Range("A205").Value = 3
$ws.Range("B205").Value = "Femacal de La Calera"
$ws.Range("C205").Value = "Coquimbo"
$ws.Range("D205").Value = 44524
$ws.Range("E205").Value = 5
$ws.Range("F205").Value = 100112017
$ws.Range("G205").Value = "Apio"
$ws.Range("H205").Value = "Americana (o)"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 160
$ws.Range("K205").Value = 9000
$ws.Range("L205").Value = 9000
$ws.Range("M205").Value = 9000
$ws.Range("N205").Value = "`$/docena de matas"
$ws.Range("O205").Value = "Pan de Azúcar"
$ws.Range("P205").Value = 1500
$ws.Range("Q205").Value = 6
$ws.Range("R205").Value = "Hortaliza"
